$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $donorAddr, $text) {
    # Force the literal string into the cell as TEXT (not an auto-parsed date)
    # by temporarily marking the cell as Text-formatted before assignment...
    $ws.Range($cellAddr).NumberFormat = "@"
    $ws.Range($cellAddr).Value = $text
    # ...then restore the original (General) cell formatting/style by pasting
    # the format from a same-row donor cell that already carries the correct
    # style, so the visible style index matches the rest of the table.
    $ws.Range($donorAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial(-4122)  # xlPasteFormats
}

# --- "Bad Drivers" block ---
$ws.Range("C3").Value = 42
$ws.Range("D3").Value = 95.3
$ws.Range("C4").Value = 42

# --- "Good Drivers" block: rows shift up one position, row 12's driver
#     wraps around to row 17 (with a refreshed client count) ---

# Row 12: was "21.60.2.1" -> now "23.100.0.4"
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B12").Value = 445055
$ws.Range("D12").Value = 99.90000000000001
Set-TextValue "E12" "D12" "2024-11-10"

# Row 13: was "22.50.1.1" -> now "22.80.0.9"
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B13").Value = 77849
$ws.Range("D13").Value = 99.90000000000001
Set-TextValue "E13" "D13" "2021-08-18"

# Row 14: was "23.100.0.4" -> now "22.50.1.1"
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
Set-TextValue "E14" "D14" "2021-04-27"

# Row 15: was "22.80.0.9" -> now "21.110.3.2"
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B15").Value = 59673
$ws.Range("D15").Value = 100
Set-TextValue "E15" "D15" "2020-08-05"

# Row 16: was "21.110.3.2" -> now "21.70.0.6"
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B16").Value = 113652
$ws.Range("D16").Value = 100
Set-TextValue "E16" "D16" "2020-01-06"

# Row 17: was "21.70.0.6" -> now "21.60.2.1" (D17/E17 stay as-is)
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B17").Value = 56018
